# Prepped for EJ46FF measurements.
# Fill in the real string-gauge headers on the EJ46FF sheet (previously
# generic "String 1".."String 6" placeholders) and make EJ46FF the
# active/selected sheet, matching the new measurement data about to be
# collected there.

$wb = $excel.ActiveWorkbook

$ej45ff = $wb.Worksheets.Item("EJ45FF")
$ej46ff = $wb.Worksheets.Item("EJ46FF")

# Replace the placeholder headers with the real EJ46FF part numbers.
$ej46ff.Range("B1").Value = "J4601FF"
$ej46ff.Range("C1").Value = "J4602FF"
$ej46ff.Range("D1").Value = "J4603FF"
$ej46ff.Range("E1").Value = "J4604FF"
$ej46ff.Range("F1").Value = "J4605FF"
$ej46ff.Range("G1").Value = "J4606FF"

# Move the selection on EJ45FF up to the header row, and clear its
# "tab selected" status in favor of EJ46FF.
$ej45ff.Range("B1:G1").Select()

# EJ46FF becomes the active/selected sheet, with the single cell B2
# selected (ready for the first new measurement entry).
$ej46ff.Activate()
$ej46ff.Range("B2").Select()
